# Shift four text boxes on slide 2 to the right (x-offset only; y unchanged).
# Commit message: "update brave img, 컨트리뷰톤 Title"
#
# Target x offsets (EMU) -> converted to points (1 pt = 12700 EMU) for the
# Shape.Left property. The literal point values below are chosen so that,
# after PowerPoint's internal Single (float32) rounding of .Left, the EMU
# value written back into the XML lands exactly on the intended target.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$moves = @(
    @{ Name = "TextBox 3"; LeftPt = 274.4620819091118 },  # off x: 3292628 -> 3485668
    @{ Name = "TextBox 5"; LeftPt = 310.10917663578914 }, # off x: 3745346 -> 3938386
    @{ Name = "TextBox 7"; LeftPt = 556.5886230469356 },  # off x: 6875635 -> 7068675
    @{ Name = "TextBox 8"; LeftPt = 317.4773559569123 }   # off x: 3838922 -> 4031962
)

foreach ($move in $moves) {
    $shp = $s.Shapes.Item($move.Name)
    $shp.Left = $move.LeftPt
}
